# Updated backlog - Angie
# Rewrites Sheet1's task list (columns B/C duplicate task names, D mostly
# cleared out, F column becomes a blank "Posted" stub column down through a
# couple of freshly-added trailing rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlVAlignCenter = -4108
$xlVAlignTop = -4160

function Set-Plain($row, $col, $text) {
    # Value with no special alignment/wrap style (style index 0 / none).
    $c = $ws.Cells.Item($row, $col)
    $c.Clear()
    $c.Value2 = $text
}

function Set-Wrapped($row, $col, $text) {
    # wrapText only -> style index 1
    $c = $ws.Cells.Item($row, $col)
    $c.Clear()
    $c.Value2 = $text
    $c.WrapText = $true
}

function Set-WrappedCenter($row, $col, $text) {
    # wrapText + vertical center -> style index 2
    $c = $ws.Cells.Item($row, $col)
    $c.Clear()
    $c.Value2 = $text
    $c.WrapText = $true
    $c.VerticalAlignment = $xlVAlignCenter
}

function Set-WrappedTop($row, $col, $text) {
    # wrapText + vertical top -> style index 3
    $c = $ws.Cells.Item($row, $col)
    $c.Clear()
    $c.Value2 = $text
    $c.WrapText = $true
    $c.VerticalAlignment = $xlVAlignTop
}

# --- Row 2 ---
Set-Wrapped 2 2 "correct sub-title on home pg (BEN)"
Set-Wrapped 2 3 "correct sub-title on home pg (BEN)"
Set-Wrapped 2 6 ""

# --- Row 3 ---
Set-WrappedCenter 3 2 "create contact pg (Angie)"
Set-WrappedCenter 3 3 "create contact pg (Angie)"
Set-Wrapped 3 4 ""
Set-WrappedTop 3 6 ""

# --- Row 4 ---
Set-WrappedCenter 4 2 "3NF tables w/ JOINS (Rebecca)"
Set-WrappedCenter 4 3 "3NF tables w/ JOINS (Rebecca)"
$ws.Cells.Item(4, 4).Clear()
Set-Wrapped 4 6 ""

# --- Row 5 ---
Set-Wrapped 5 2 "Confirm format (all)"
Set-Wrapped 5 3 "Confirm format (all)"
Set-Wrapped 5 4 ""
Set-Wrapped 5 6 ""

# --- Row 6 ---
Set-Wrapped 6 2 "Query scout with most sales (Jim)"
Set-Wrapped 6 3 "Query scout with most sales (Jim)"
$ws.Cells.Item(6, 4).Clear()
Set-Wrapped 6 6 ""

# --- Row 7 ---
Set-Plain 7 2 "Testing(Angie)"
Set-Plain 7 3 "Testing(Angie)"
$ws.Cells.Item(7, 4).Clear()
Set-Wrapped 7 6 ""

# --- Row 8 ---
Set-Wrapped 8 2 "Update price field (Angie)"
Set-Wrapped 8 3 "Update price field (Angie)"

# --- Row 9 ---
Set-Wrapped 9 2 "Rewrite queries(all)"
Set-Wrapped 9 3 "Rewrite queries(all)"
$ws.Cells.Item(9, 4).Clear()
Set-Wrapped 9 6 ""

# --- Row 10 (blank stub cells) ---
Set-Wrapped 10 3 ""
Set-Wrapped 10 4 ""
Set-Wrapped 10 5 ""
Set-Wrapped 10 6 ""

# --- Row 11 (blank stub cells) ---
Set-Wrapped 11 3 ""
Set-Wrapped 11 6 ""

# --- Rows 12-13 (new trailing rows, blank stub cells) ---
Set-Wrapped 12 6 ""
Set-Wrapped 13 6 ""

# --- sheet view: scrolled down a bit with a different active selection ---
$ws.Application.Goto($ws.Range("A3"))
$ws.Range("A9").Select()
